$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 78 - this shifts the existing
# rows 78-88 down to 79-89 and keeps all their data/formatting intact.
$ws.Rows.Item(78).Insert()

# Populate the newly inserted row 78 with the new weekly record.
$ws.Range("A78").Value = 10
$ws.Range("B78").Value = "Vega Modelo de Temuco"
$ws.Range("C78").Value = "La Araucanía"
$ws.Range("D78").Value = 44776
$ws.Range("E78").Value = 9
$ws.Range("F78").Value = "Fruta"
$ws.Range("G78").Value = 100108
$ws.Range("H78").Value = "Tropicales y subtropicales"
$ws.Range("I78").Value = 100108007
$ws.Range("J78").Value = "Coco"
$ws.Range("K78").Value = "Sin especificar"
$ws.Range("L78").Value = "Primera"
$ws.Range("M78").Value = 40
$ws.Range("N78").Value = 30000
$ws.Range("O78").Value = 30000
$ws.Range("P78").Value = 30000
$ws.Range("Q78").Value = "$/malla 20 unidades"
$ws.Range("R78").Value = "Perú"
$ws.Range("S78").Value = 1500
$ws.Range("T78").Value = 20
